$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New pretty-printed JSON-like text replacing the old single-line dict text
$newText = @'
questions = [
    {
        "title": "You work for a retail company that wishes to migrate its on-premises transactional data to Azure. You have been tasked with the development of a managed, cloud-based data integration service that can ingest, prepare, transform, and transfer this data to Azure storage solutions at scale.Which Azure service should you use?",
        "ques_type": 2,
        "options": [
            "Azure Blob Storage",
            "Azure Data Lake",
            "Azure Data Factory",
            "Azure Cosmos DB"
        ],
        "score": "Azure Data Factory"
    },
    {
        "title": "You work for a financial institution that has raw transactional data that requires cleansing, transformation, and enrichment before it can be used for analytics. They need an Azure service that integrates seamlessly with Azure Data Lake and provides a rich set of transformations. You have been tasked with setting up the required infrastructure.Which Azure service should you use?",
        "ques_type": 2,
        "options": [
            "Azure Stream Analytics",
            "Azure Databricks",
            "Azure Logic Apps",
            "Azure Synapse Analytics"
        ],
        "score": "Azure Databricks"
    },
    {
        "title": "You work for a multinational retail chain that aims to derive insights from its vast customer purchase data. The company is already using Azure for its infrastructure, and you have been tasked with visually representing trends and anomalies to make data-driven decisions. What should you do?",
        "ques_type": 2,
        "options": [
            "Execute a series of T-SQL queries to generate structured reports.",
            "Utilize Power Query to extract and transform data into visual models.",
            "Implement a REST API to fetch and display data in graphical form.",
            "Convert raw data to CSV and analyze using Azure Databricks."
        ],
        "score": "Utilize Power Query to extract and transform data into visual models."
    },
    {
        "title": "You are a data scientist at an online travel agency. The agency has a vast database of hotels and destinations on Azure SQL Database, which is experiencing slow query performance. You have been tasked with optimizing the performance of the queries.What should you do to accomplish this with the least amount of effort?",
        "ques_type": 2,
        "options": [
            "Partition the data tables based on high-frequency access patterns.",
            "Implement Automatic Tuning to continuously adapt to changing workloads.",
            "Introduce indexing on frequently queried columns to improve read performance.",
            "Increase the storage size of the Azure SQL Database."
        ],
        "score": "Implement Automatic Tuning to continuously adapt to changing workloads."
    }
]
'@

# The original sheet had two rows: A1 = 0 (bold, bordered, centered style) and
# A2 = the long descriptive text (default style). The target keeps only a
# single row: A1 = the (now reformatted) text, with the default/no style.
# Remove the old row 1 entirely so the old A2 (already default-styled)
# shifts up to become the new A1.
$ws.Rows.Item(1).Delete()

# Write the updated text into the now-single A1 cell. Its formatting already
# carried over as the default (unstyled) format from the old A2.
$ws.Range("A1").Value = $newText

# Writing a multi-line value auto-expands the row height; restore the
# sheet's normal auto height so row 1 keeps the default (non-custom) height.
$ws.Rows.Item(1).AutoFit()
